$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that needs to move
# from 45334 (2024-02-12) to 45335 (2024-02-13) for rows 2 through 27.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = 45335
}
